# TC08_Canine_Filter_StageOfDisease-5.xlsx - "10 icdc scripts for jenkins"
#
# The FilesTab Neo4j query (row 4, column B on the "startup" sheet) is
# updated to drop the `File Type` and `Breed` columns from its RETURN
# clause.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newFilesQuery = @'

MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE diag.stage_of_disease IN ['V']
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $newFilesQuery

# Match the author's new selection/scroll position on the worksheet.
$ws.Activate()
$ws.Range("B4").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
